$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row: only J1 ("Editorial") is cleared out; its highlighted
# header style is downgraded to the plain style (matching A9, which is
# already style s="1") by copying formats only.
# ---------------------------------------------------------------------
$ws.Range("J1").Value = ""
$ws.Range("A9").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Data rows 2-6: content for 2nd Sept is corrected / re-pointed.
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "https://leetcode.com/problems/3sum/description/"

$ws.Range("C3").Value = "Two-Pointer"
$ws.Range("D3").Value = "Medium"
$ws.Range("E3").Value = " 3Sum"
$ws.Range("F3").Value = "https://leetcode.com/problems/3sum/description/"

# Rows 4-6 keep identical text (only shared-string indices shifted in the
# source diff), so nothing else to change here.

# ---------------------------------------------------------------------
# Data rows 7-12: newly filled in with the "2nd Sept" entries.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Array"
$ws.Range("C7").Value = "Matrix"
$ws.Range("D7").Value = "Medium"
$ws.Range("E7").Value = "Rotate Image"
$ws.Range("F7").Value = "https://leetcode.com/problems/rotate-image/description/"
$ws.Range("G7").Value = 45537
$ws.Range("H7").Value = "Sept"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Array"
$ws.Range("C8").Value = "Logic"
$ws.Range("D8").Value = "Medium"
$ws.Range("E8").Value = " String to Integer (atoi)"
$ws.Range("F8").Value = "https://leetcode.com/problems/string-to-integer-atoi/description/"
$ws.Range("G8").Value = 45537
$ws.Range("H8").Value = "Sept"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Array"
$ws.Range("C9").Value = "Logic"
$ws.Range("D9").Value = "Easy"
$ws.Range("E9").Value = " Integer to Roman"
$ws.Range("F9").Value = "https://leetcode.com/problems/roman-to-integer/description/"
$ws.Range("G9").Value = 45537
$ws.Range("H9").Value = "Sept"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Array"
$ws.Range("C10").Value = "Hash Table"
$ws.Range("D10").Value = "Easy"
$ws.Range("E10").Value = "Roman to Integer"
$ws.Range("F10").Value = "https://leetcode.com/problems/roman-to-integer/description/"
$ws.Range("G10").Value = 45537
$ws.Range("H10").Value = "Sept"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "String"
$ws.Range("C11").Value = "Two-Pointer, String Matching"
$ws.Range("D11").Value = "Easy"
$ws.Range("E11").Value = " Implement strStr()"
$ws.Range("F11").Value = "https://leetcode.com/problems/find-the-index-of-the-first-occurrence-in-a-string/description/"
$ws.Range("G11").Value = 45537
$ws.Range("H11").Value = "Sept"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "String"
$ws.Range("C12").Value = "Sliding Window, Hash Table"
$ws.Range("D12").Value = "Medium"
$ws.Range("E12").Value = "Longest Substring W/O Repeating Characters"
$ws.Range("F12").Value = "https://leetcode.com/problems/longest-substring-without-repeating-characters/description/"
$ws.Range("G12").Value = 45537
$ws.Range("H12").Value = "Sept"

# ---------------------------------------------------------------------
# Fix up number formats / styles that need to change now that rows 7-12
# carry real data: G9:G11 need the "d-mmm" date style (matching G2), and
# H7:H12 need the "mmm-yy" month style (matching H2). Rows 7,8,12 already
# had the G/H styles pre-set, so re-applying is harmless.
# ---------------------------------------------------------------------
$ws.Range("G2").Copy()
$ws.Range("G9:G11").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("H7:H12").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Restore the selection the author left the sheet with.
# ---------------------------------------------------------------------
$ws.Range("E14").Select()
